$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 4, shifting existing rows 4..45 down to 5..46
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new entry
$ws.Cells.Item(4, 1).Value = "ESCAPE CLEANING COMPANY LLC"
$ws.Cells.Item(4, 2).Value = "Zigan, Gerald L"
$ws.Cells.Item(4, 3).Value = "030"
$ws.Cells.Item(4, 4).Value = Get-Date -Year 2025 -Month 9 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(4, 5).Value = "0008145"
